$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 12.7228780040422

$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 2.461472421528573
